# Update cryptos list figures (prices / volume deltas / reordered rows)
# as published by the scheduled GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (e.g. "66.993.18", "1.00")
# -- mark the cells we are about to rewrite as Text first so Excel
# does not auto-coerce the numeric-looking strings into numbers.
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D22","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D34","D36","D39","D40","D41","D42","D43","D46","D47","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.395.71"
$ws.Range("E2").Value = "  +5.07%  "
$ws.Range("D3").Value = "3.462.26"
$ws.Range("E3").Value = "  +4.28%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "581.64"
$ws.Range("E5").Value = "  +5.28%  "
$ws.Range("D6").Value = "185.22"
$ws.Range("E6").Value = "  +7.45%  "
$ws.Range("D7").Value = "0.632"
$ws.Range("E7").Value = "  +2.50%  "
$ws.Range("D8").Value = "3.458.07"
$ws.Range("E8").Value = "  +4.39%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +1.46%  "
$ws.Range("D11").Value = "0.647"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "56.58"
$ws.Range("E12").Value = "  +5.75%  "
$ws.Range("D13").Value = "0.0000279"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "9.45"
$ws.Range("E14").Value = "  +4.62%  "
$ws.Range("D15").Value = "3.999.08"
$ws.Range("E15").Value = "  +4.03%  "
$ws.Range("D16").Value = "18.70"
$ws.Range("E16").Value = "  +3.41%  "
$ws.Range("D17").Value = "3.455.25"
$ws.Range("E17").Value = "  +4.53%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "67.239.82"
$ws.Range("E18").Value = "  +4.88%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "0.120"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "12.13"
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("D22").Value = "481.35"
$ws.Range("E22").Value = "  +5.65%  "
$ws.Range("D23").Value = "5.43"
$ws.Range("E23").Value = "  +8.31%  "
$ws.Range("D24").Value = "16.78"
$ws.Range("E24").Value = "  +21.39%  "
$ws.Range("D25").Value = "4.43"
$ws.Range("E25").Value = "  +8.81%  "
$ws.Range("D26").Value = "89.69"
$ws.Range("E26").Value = "  +2.87%  "
$ws.Range("D27").Value = "11.04"
$ws.Range("E27").Value = "  +3.40%  "
$ws.Range("E28").Value = "  +3.09%  "
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  +7.01%  "
$ws.Range("D30").Value = "31.44"
$ws.Range("E30").Value = "  +2.21%  "
$ws.Range("D31").Value = "7.15"
$ws.Range("E31").Value = "  +9.46%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "64.39"
$ws.Range("E32").Value = "  +5.76%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "11.73"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("D34").Value = "591.61"
$ws.Range("E34").Value = "  +4.59%  "
$ws.Range("E35").Value = "  +5.34%  "
$ws.Range("D36").Value = "0.149"
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("D39").Value = "36.59"
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0773"
$ws.Range("E40").Value = "  +6.22%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "0.386"
$ws.Range("E41").Value = "  +5.64%  "
$ws.Range("D42").Value = "3.206.24"
$ws.Range("E42").Value = "  +5.05%  "
$ws.Range("D43").Value = "2.92"
$ws.Range("E43").Value = "  +6.11%  "
$ws.Range("E44").Value = "  +4.26%  "
$ws.Range("E45").Value = "  +4.41%  "
$ws.Range("D46").Value = "2.79"
$ws.Range("E46").Value = "  +22.97%  "
$ws.Range("D47").Value = "3.21"
$ws.Range("E47").Value = "  +2.01%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "8.78"
$ws.Range("E49").Value = "  +8.38%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").Value = "0.997"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("E51").Value = "  +9.37%  "
